$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.3631899612946223
$ws.Range("C2").Value = 3.8824983095132808
$ws.Range("D2").Value = 3.6330016612950375
$ws.Range("E2").Value = 0.8255980946300262

$ws.Range("B3").Value = 1.3066119353810981
$ws.Range("C3").Value = 5.271465140788357
$ws.Range("D3").Value = 1.6186022660001569
$ws.Range("E3").Value = 5.0968846257258793

$ws.Range("B4").Value = 2.744953431473101
$ws.Range("C4").Value = 2.9953444281333996
$ws.Range("D4").Value = 5.3167672498780094
$ws.Range("E4").Value = 3.9005271121602982

$ws.Range("B5").Value = 1.1878556071561832
$ws.Range("C5").Value = 1.2007418807514199
$ws.Range("D5").Value = 1.9279333932121985
$ws.Range("E5").Value = 1.1867814052409109

$ws.Range("E6").Value = 0.75087811552978445

$ws.Range("B9").Value = 121.3979561412464
$ws.Range("C9").Value = 103.65609784977876
$ws.Range("D9").Value = 117.09524216621371
$ws.Range("E9").Value = 126.54394667475114

$ws.Range("E10").Value = 25.601111440288133

$ws.Range("B11").Value = 1.3999999994566925
$ws.Range("C11").Value = 1.3999998355699337
$ws.Range("D11").Value = 1.3999786774595075
$ws.Range("E11").Value = 1.3999979751247735

$ws.Range("B12").Value = 2.7649400953948122
$ws.Range("C12").Value = 4.6617607344072516
$ws.Range("D12").Value = 5.4138478873375
$ws.Range("E12").Value = 0.96754328040951476

$ws.Range("B13").Value = 11.392702095153142
$ws.Range("C13").Value = 13.08895253319049
$ws.Range("D13").Value = 8.7790992131708077
$ws.Range("E13").Value = 1.5023597780173634

$ws.Range("B14").Value = 2.5452977567500756
$ws.Range("C14").Value = 4.1002052642163722
$ws.Range("D14").Value = 1.7025443527128425
$ws.Range("E14").Value = 1.8395614034637076

$ws.Range("B16").Value = 1.0254529775675008
$ws.Range("C16").Value = 1.0410020526421637
$ws.Range("D16").Value = 1.0170254435271173
$ws.Range("E16").Value = 1.0183956140346371

$ws.Range("B17").Value = 9.8007749295388731
$ws.Range("C17").Value = 8.545750363861643
$ws.Range("D17").Value = 14.164002362755896
$ws.Range("E17").Value = 1.9197927718678569

$ws.Range("B18").Value = 5.9094172135756997
$ws.Range("C18").Value = 3.7710212715871574
$ws.Range("D18").Value = 3.6890425452043223
$ws.Range("E18").Value = 11.12240557975192

$ws.Range("B19").Value = 68.578349895128355
$ws.Range("C19").Value = 64.856945603168015
$ws.Range("D19").Value = 72.110912306016473
$ws.Range("E19").Value = 68.905032887965518
